# Vaccine workbook fixes:
#  - split the combined "10 pack...vials 5 pack...syringes" Kinrix packaging
#    string into two distinct rows
#  - correct several Vaccine-name typos / missing spaces / slashes
#  - normalize spacing on a few Influenza / Tetanus / Pneumococcal labels

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Pediatric VFC Vaccine ---
$ws1 = $wb.Worksheets.Item("Pediatric VFC Vaccine ")

# Kinrix (DTaP-IPV): split the combined packaging text across the two rows
$ws1.Range("D6").Value = "10 pack - 1 dose vials "
$ws1.Range("D7").Value = "5 pack - 1 dose T-L syringes "

# TriHIBit row: Vaccine column was "DTaP ", should be "DTaP-Hib "
$ws1.Range("A11").Value = "DTaP-Hib "

# COMVAX row: Vaccine column was "Hepatitis B ", should be "Hepatitis B-Hib "
$ws1.Range("A14").Value = "Hepatitis B-Hib "

# ENGERIX B rows: missing space/slash
$ws1.Range("A20").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A21").Value = "Hepatitis B Pediatric/Adolescent"

# RECOMBIVAX HB row: missing slash (keeps its existing double space before "Pediatric")
$ws1.Range("A22").Value = "Hepatitis B  Pediatric/Adolescent"

# Prevnar row: missing space
$ws1.Range("A30").Value = "Pneumococcal 7-valent (Pediatric)"

# Pneumovax row: missing space
$ws1.Range("A31").Value = "Pneumococcal Polysaccharide (23 Valent)"

# --- Sheet 2: Adult VFC Vaccine ---
$ws2 = $wb.Worksheets.Item("Adult VFC Vaccine ")

# Tetanus Diphtheria Toxoids (Adsorbed for Adults) rows: missing space
$ws2.Range("A12").Value = "Tetanus  Diphtheria Toxoids "
$ws2.Range("A13").Value = "Tetanus  Diphtheria Toxoids "

# --- Sheet 3: Pediatric influenza Influenza ---
$ws3 = $wb.Worksheets.Item("Pediatric influenza Influenza")

$ws3.Range("A2").Value = "Influenza   (Age 6 months and older)"
$ws3.Range("A3").Value = "Influenza  (Age 6-35 months)"
$ws3.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("A8").Value = "Influenza  Live, Intranasal (Age 2-49 years)"
